$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SolverSettings")

$ws.Range("A10").Value = "include_RPS"
$ws.Range("B10").Value = "Y"
$ws.Range("C10:AA10").Value = "N"

$ws.Activate()
